$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New food entries to append below the existing table (rows 17..43).
# Columns: A=food, B=salty, C=effort, D=takeaway
$rows = @(
    @('Mac and Cheese', 'herzhaft', 'mittel', 'kochen'),
    @('Wraps', 'herzhaft', 'mittel', 'kochen'),
    @('Instant Nudeln', 'herzhaft', 'wenig', 'kochen'),
    @('Pilz-Risotto', 'herzhaft', 'hoch', 'kochen'),
    @('Pfannekuchen', 'süß', 'mittel', 'kochen'),
    @('Milchreis', 'süß', 'hoch', 'kochen'),
    @('Bowl', 'herzhaft', 'bestellen', 'bestellen'),
    @('Sommerrollen', 'herzhaft', 'bestellen', 'bestellen'),
    @('Burger', 'herzhaft', 'bestellen', 'bestellen'),
    @('Miracoli', 'herzhaft', 'wenig', 'kochen'),
    @('Muscheln in Weißweinsoße', 'herzhaft', 'hoch', 'kochen'),
    @('Lasagne', 'herzhaft', 'hoch', 'kochen'),
    @('Soja-Bolognese', 'herzhaft', 'mittel', 'kochen'),
    @('Kartoffeln mit Quark', 'herzhaft', 'mittel', 'kochen'),
    @('Griesbrei mit Apfelmus', 'süß', 'wenig', 'kochen'),
    @('Döner', 'herzhaft', 'bestellen', 'bestellen'),
    @('Paneer Butter Masala', 'herzhaft', 'bestellen', 'bestellen'),
    @('Pho Suppe', 'herzhaft', 'bestellen', 'bestellen'),
    @('Asia-Nudeln', 'herzhaft', 'bestellen', 'bestellen'),
    @('Donut', 'süß', 'bestellen', 'bestellen'),
    @('Käsespätzle', 'herzhaft', 'hoch', 'kochen'),
    @('Nudeln mit Pesto', 'herzhaft', 'wenig', 'kochen'),
    @('Gnocchi', 'herzhaft', 'wenig', 'kochen'),
    @('Maultauschen', 'herzhaft', 'wenig', 'kochen'),
    @('Sandwiches', 'herzhaft', 'wenig', 'kochen'),
    @('Ikea Köttbullar', 'herzhaft', 'bestellen', 'bestellen'),
    @('Beck Fladen', 'herzhaft', 'bestellen', 'bestellen')
)

$startRow = 17
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
}

$win = $excel.ActiveWindow
$win.ScrollRow = 35
$ws.Range("E40").Select()
